$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 348
$wsExpo.Range("F4").Value = 2916
$wsExpo.Range("F5").Value = 72
$wsExpo.Range("F6").Value = 614

# Sheet "全部类型" (All Types) - update "想去人数" (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 348
$wsAll.Range("F6").Value = 2916
$wsAll.Range("F7").Value = 72
$wsAll.Range("F8").Value = 614
